$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct text/string updates (values Excel will not misinterpret as numbers)
$ws.Range('D2').Value = '29.863.38'
$ws.Range('D3').Value = '1.886.60'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('E5').Value = '  -5.08%  '
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  -1.42%  '
$ws.Range('E9').Value = '  -2.05%  '
$ws.Range('E10').Value = '  -2.69%  '
$ws.Range('E11').Value = '  +4.94%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.972.75'
$ws.Range('E12').Value = '  +4.26%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('E14').Value = '  -2.87%  '
$ws.Range('E15').Value = '  -1.06%  '
$ws.Range('E16').Value = '  -1.88%  '
$ws.Range('D17').Value = '29.927.80'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('E19').Value = '  -1.58%  '
$ws.Range('E20').Value = '  -0.76%  '
$ws.Range('D21').Value = '2.148.33'
$ws.Range('E21').Value = '  +2.18%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  -2.37%  '
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('E26').Value = '  -1.50%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('E30').Value = '  +5.12%  '
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('E32').Value = '  -0.57%  '
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('E34').Value = '  -4.03%  '
$ws.Range('E35').Value = '  -1.43%  '
$ws.Range('E36').Value = '  -1.82%  '
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('E38').Value = '  +1.31%  '
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  -0.64%  '
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('D42').Value = '1.101.41'
$ws.Range('E42').Value = '  -4.15%  '
$ws.Range('E43').Value = '  +1.40%  '
$ws.Range('E44').Value = '  -2.43%  '
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('E50').Value = '  -2.98%  '
$ws.Range('D51').Value = '2.056.30'
$ws.Range('E51').Value = '  +2.07%  '

# Numeric-looking text values must be forced to remain text (matching the original
# inlineStr type) instead of being auto-converted by Excel into real numbers.
# For each one: write a formula that evaluates to the exact text, then collapse it
# to a static value via Copy + PasteSpecial (values only), one cell at a time so the
# paste lands on the correct cell.
$ws.Range('D4').Formula = '="1.002"'
$ws.Range('D4').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('D5').Formula = '="0.7479"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('D6').Formula = '="242.43"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('D7').Formula = '="1.001"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('D8').Formula = '="0.3117"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('D9').Formula = '="25.34"'
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('D10').Formula = '="0.07122"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('D11').Formula = '="0.08523"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('D13').Formula = '="0.7595"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('D14').Formula = '="5.357"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('D15').Formula = '="93.41"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('D16').Formula = '="6.146"'
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('D19').Formula = '="242.98"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('D20').Formula = '="0.000007788"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('D23').Formula = '="7.982"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('D24').Formula = '="1.002"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('D25').Formula = '="0.1594"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('D26').Formula = '="9.356"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('D27').Formula = '="162.71"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('D28').Formula = '="18.74"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('D29').Formula = '="2.026"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('D30').Formula = '="1.515"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('D31').Formula = '="1.531"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('D32').Formula = '="4.473"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('D33').Formula = '="4.096"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('D34').Formula = '="0.05394"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('D36').Formula = '="0.7426"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('D37').Formula = '="1.003"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('D39').Formula = '="0.01937"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('D40').Formula = '="2.772"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('D41').Formula = '="0.4449"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('D44').Formula = '="72.38"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('D45').Formula = '="0.8558"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('D47').Formula = '="102.52"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('D48').Formula = '="7.658"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('D49').Formula = '="1.861"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)

$excel.CutCopyMode = 0

Write-Output "Updated cryptos list"